$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44446
$ws.Range("L2").Value = "Primera"
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("S2").Value = 2150

# Row 3
$ws.Range("D3").Value = 44487
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 23500
$ws.Range("S3").Value = 2350

# Row 4
$ws.Range("D4").Value = 44447

# Row 5
$ws.Range("D5").Value = 44461
$ws.Range("L5").Value = "Especial"
$ws.Range("N5").Value = 31000
$ws.Range("O5").Value = 32000
$ws.Range("P5").Value = 31500
$ws.Range("S5").Value = 3150

# Row 6
$ws.Range("D6").Value = 44461
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("S6").Value = 3000

# Row 9
$ws.Range("D9").Value = 44448

# Row 10
$ws.Range("D10").Value = 44452
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 21000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21500
$ws.Range("S10").Value = 2150
